$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The two "Resolving-Mac" rows (old rows 4 and 5) are dropped entirely,
# which also drops the last reference to the "Resolving-Mac" shared string.
$ws.Rows.Item(4).Delete()
$ws.Rows.Item(4).Delete()

# Update row 2
$ws.Range("B2").Value = "Rbp4"
$ws.Range("C2").Value = "Stra6"
$ws.Range("G2").Value = 4.164095000000001
$ws.Range("I2").Value = 1
$ws.Range("J2").Value = 1
$ws.Range("O2").Value = 0.6332315371308455
$ws.Range("P2").Value = 0.6332315371308455
$ws.Range("Q2").Value = 1.431289673558334
$ws.Range("S2").Value = 0.6332315371308455
$ws.Range("T2").Value = 0.6332315371308455

# Update row 3
$ws.Range("B3").Value = "Rbp4"
$ws.Range("C3").Value = "Stra6"
$ws.Range("D3").Value = "MuSCs"
$ws.Range("G3").Value = 4.164095000000001
$ws.Range("I3").Value = 1
$ws.Range("J3").Value = 1
$ws.Range("M3").Value = 0.199084
$ws.Range("N3").Value = 0.597252
$ws.Range("O3").Value = 0.3667684628691545
$ws.Range("P3").Value = 0.3667684628691545
$ws.Range("Q3").Value = 0.8290046889800001
$ws.Range("R3").Value = 7.461042200820001
$ws.Range("S3").Value = 0.3667684628691545
$ws.Range("T3").Value = 0.3667684628691545
